$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete the rows for "CAR" (row 2) and "GAE" (row 4).
# Delete row 4 first so row indices for the remaining deletion stay valid.
$ws.Rows.Item(4).Select()
$ws.Rows.Item(4).Delete()

$ws.Rows.Item(2).Select()
$ws.Rows.Item(2).Delete()

$ws.Range("A3:XFD3").Select()
